$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰" + [char]10 + `
"✅ Dólar paralelo: 68" + [char]10 + `
"" + [char]10 + `
"Binance" + [char]10 + `
"✅ 1000 Bs = 12.91 = 52078.76 pesos" + [char]10 + `
"✅ 52078.76 pesos = 12.82 = 961.39 Bs" + [char]10 + `
"" + [char]10 + `
"Promedio competencia" + [char]10 + `
"✅ Tasa pesos: 20" + [char]10 + `
"✅ Tasa Bs: 20" + [char]10 + `
"✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 77.479
$ws2.Range("O10").Value = 4035.01
$ws2.Range("N12").Value = 4062.79
